$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated symbol list (cryptocurrency prices/volumes/hour marker) pulled by the
# scheduled GitHub Actions sync job. Each update targets Price (D), Volume(1h)
# (E, text label - only changed for a couple of rows) and Hora (G) columns for
# rows 2-51.
$updates = @(
    @{Cell="D2"; Value="245.85"},
    @{Cell="G2"; Value="2"},
    @{Cell="D3"; Value="21.89"},
    @{Cell="G3"; Value="2"},
    @{Cell="D4"; Value="5.412"},
    @{Cell="G4"; Value="2"},
    @{Cell="D5"; Value="0.05754"},
    @{Cell="G5"; Value="2"},
    @{Cell="D6"; Value="3.403"},
    @{Cell="G6"; Value="2"},
    @{Cell="D7"; Value="6.328"},
    @{Cell="G7"; Value="2"},
    @{Cell="D8"; Value="0.8172"},
    @{Cell="G8"; Value="2"},
    @{Cell="D9"; Value="0.9852"},
    @{Cell="E9"; Value="8FTXTokenFTTBestin24h"},
    @{Cell="G9"; Value="2"},
    @{Cell="D10"; Value="0.1430"},
    @{Cell="G10"; Value="2"},
    @{Cell="D11"; Value="0.07352"},
    @{Cell="G11"; Value="2"},
    @{Cell="D12"; Value="0.03129"},
    @{Cell="G12"; Value="2"},
    @{Cell="D13"; Value="0.03019"},
    @{Cell="G13"; Value="2"},
    @{Cell="D14"; Value="4.164"},
    @{Cell="G14"; Value="2"},
    @{Cell="D15"; Value="0.09401"},
    @{Cell="G15"; Value="2"},
    @{Cell="D16"; Value="0.001585"},
    @{Cell="G16"; Value="2"},
    @{Cell="D17"; Value="0.04803"},
    @{Cell="G17"; Value="2"},
    @{Cell="D18"; Value="0.0005849"},
    @{Cell="G18"; Value="2"},
    @{Cell="D19"; Value="0.006212"},
    @{Cell="G19"; Value="2"},
    @{Cell="D20"; Value="0.004112"},
    @{Cell="G20"; Value="2"},
    @{Cell="D21"; Value="0.0009964"},
    @{Cell="G21"; Value="2"},
    @{Cell="D22"; Value="0.0001500"},
    @{Cell="G22"; Value="2"},
    @{Cell="D23"; Value="3.755"},
    @{Cell="G23"; Value="2"},
    @{Cell="D24"; Value="2.205"},
    @{Cell="G24"; Value="2"},
    @{Cell="G25"; Value="2"},
    @{Cell="G26"; Value="2"},
    @{Cell="D27"; Value="0.0003998"},
    @{Cell="G27"; Value="2"},
    @{Cell="G28"; Value="2"},
    @{Cell="G29"; Value="2"},
    @{Cell="G30"; Value="2"},
    @{Cell="G31"; Value="2"},
    @{Cell="G32"; Value="2"},
    @{Cell="G33"; Value="2"},
    @{Cell="G34"; Value="2"},
    @{Cell="G35"; Value="2"},
    @{Cell="G36"; Value="2"},
    @{Cell="G37"; Value="2"},
    @{Cell="G38"; Value="2"},
    @{Cell="G39"; Value="2"},
    @{Cell="D40"; Value="0.03886"},
    @{Cell="G40"; Value="2"},
    @{Cell="D41"; Value="0.006396"},
    @{Cell="G41"; Value="2"},
    @{Cell="D42"; Value="0.1073"},
    @{Cell="G42"; Value="2"},
    @{Cell="D43"; Value="0.002631"},
    @{Cell="G43"; Value="2"},
    @{Cell="D44"; Value="0.006680"},
    @{Cell="G44"; Value="2"},
    @{Cell="D45"; Value="0.00005611"},
    @{Cell="G45"; Value="2"},
    @{Cell="G46"; Value="2"},
    @{Cell="D47"; Value="0.3799"},
    @{Cell="E47"; Value="46CoinbaseStockTokenCOIN"},
    @{Cell="G47"; Value="2"},
    @{Cell="G48"; Value="2"},
    @{Cell="D49"; Value="0.00002100"},
    @{Cell="G49"; Value="2"},
    @{Cell="D50"; Value="0.01010"},
    @{Cell="G50"; Value="2"},
    @{Cell="G51"; Value="2"}
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    # D/G columns carry numeric-looking text (e.g. "245.85", "2"). The source
    # workbook stores these as plain text (inline strings), so force Text
    # format before assigning the value to avoid Excel auto-converting the
    # numeric-looking string into a real number.
    $col = $u.Cell.Substring(0, 1)
    if ($col -eq "D" -or $col -eq "G") {
        $cell.NumberFormat = "@"
    }
    $cell.Value = $u.Value
}

Write-Host "Updated symbol list on Fri Dec 23 02:07:21 UTC 2022 with GitHub Actions"
